$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 - "Build API for upload": Actual time updated, comment expanded
$ws.Range("E10").Value = "45 minutes"
$ws.Range("F10").Value = "For creating and testing using Postman, Made changes in uploading json from file url and also using chunks"

# Row 11 - content shifted to "Build API for similarity check"
$ws.Range("A11").Value = "Build API for similarity check"
$ws.Range("B11").Value = "Similarity search endpoint to query documents"
$ws.Range("E11").Value = "45 minutes"
$ws.Range("F11").Value = "For creating and testing using Postman, Made changes in the response which took some time to debug"

# Row 12 - content shifted to "Build API for get"
$ws.Range("A12").Value = "Build API for get"
$ws.Range("B12").Value = "endpoint to get details of specific journal when requested"
$ws.Range("E12").Value = "30 minutes"
$ws.Range("F12").Value = "'"
